# Draftbook enhancement: refresh the four DFS "stack" tables on Sheet1 with
# a new day's slate (Milwaukee Brewers lefties / Colorado Rockies lefties /
# Texas Rangers hitters / Chicago White Sox hitters) and clear out the
# now-empty second "failure" block on the Rockies table (columns F:I of
# row 16) that used to total the righty-pitcher numbers.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---- Table 1 (A1:D7) -> Milwaukee Brewers lefties (FD, DK) ----
$ws1.Range("A1").Value = "Milwaukee Brewers lefties (FD, DK)"

$ws1.Range("A3").Value = "Yelich"
$ws1.Range("B3").Value = 3300
$ws1.Range("C3").Value = 15.4

$ws1.Range("A4").Value = "Shaw"
$ws1.Range("B4").Value = 3700
$ws1.Range("C4").Value = 12.7

$ws1.Range("A5").Value = "Villar"
$ws1.Range("B5").Value = 2400
$ws1.Range("C5").Value = 9.5

$ws1.Range("A6").ClearContents()
$ws1.Range("B6").ClearContents()
$ws1.Range("C6").ClearContents()

# ---- Table 2 (F1:I7) -> Colorado Rockies lefties (DK) ----
$ws1.Range("F1").Value = "Colorado Rockies lefties (DK)"

$ws1.Range("F3").Value = "Blackmon"
$ws1.Range("G3").Value = 5600
$ws1.Range("H3").Value = 10

$ws1.Range("F4").Value = "Dahl"
$ws1.Range("G4").Value = 4000
$ws1.Range("H4").Value = 14

$ws1.Range("F5").Value = "Gonzalez"
$ws1.Range("G5").Value = 3600
$ws1.Range("H5").Value = 27

$ws1.Range("F6").ClearContents()
$ws1.Range("G6").ClearContents()
$ws1.Range("H6").ClearContents()

# ---- Table 3 (K1:N7) -> Texas Rangers hitters (DK) ----
$ws1.Range("K1").Value = "Texas Rangers hitters (DK)"

$ws1.Range("K3").Value = "DeShields"
$ws1.Range("L3").Value = 3300
$ws1.Range("M3").Value = 10

$ws1.Range("K4").Value = "Choo"
$ws1.Range("L4").Value = 4100
$ws1.Range("M4").Value = 7

$ws1.Range("K5").Value = "Kiner-Falefa"
$ws1.Range("L5").Value = 3500
$ws1.Range("M5").Value = 3

$ws1.Range("K6").Value = "Gallo"
$ws1.Range("L6").Value = 3800
$ws1.Range("M6").Value = 14

# ---- Table 4 (A10:D16) -> Chicago White Sox hitters (FD, DK) ----
$ws1.Range("A10").Value = "Chicago White Sox hitters (FD, DK)"

$ws1.Range("A12").Value = "Anderson"
$ws1.Range("B12").Value = 3000
$ws1.Range("C12").Value = 6.2

$ws1.Range("A13").Value = "Abreu"
$ws1.Range("B13").Value = 3800
$ws1.Range("C13").Value = 0

$ws1.Range("A14").Value = "Rondon"
$ws1.Range("B14").Value = 2900
$ws1.Range("C14").Value = 18.7

$ws1.Range("A15").Value = "Palka"
$ws1.Range("B15").Value = 2900
$ws1.Range("C15").Value = 12

# The second (now-unused) mini-table on the Chicago block (F16:I16 summary
# of G12:G15/H12:H15) no longer applies to this slate - blank it out.
$ws1.Range("F16").ClearContents()
$ws1.Range("G16").ClearContents()
$ws1.Range("H16").ClearContents()

# Stray leftover label a few rows below the sheet's used range.
$ws1.Range("E21").Value = "    "

# ---- Selections, matching the saved UI state ----
$ws1.Range("M8").Select()
$ws2.Range("A1:A4").Select()

$ws1.Activate()
